# Making mercenaries cheaper to hire; Display cash on all shops; Lots of tweaks
#
# Adds a new "Mercenary Cost" sheet (modelled on the existing
# "Levels and Experience" sheet) that computes a per-level mercenary
# hiring cost, and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Reference sheet used as a styling template (same B/C column layout).
$wsLevels = $wb.Worksheets.Item("Levels and Experience")

# --- Update "Levels and Experience" selection -----------------------------
[void]$wsLevels.Range("A1:H27").Select()

# --- Create the new "Mercenary Cost" sheet, after the last sheet ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Mercenary Cost"

# Headers
$ws.Range("B2").Value = "Level"
$ws.Range("C2").Value = "Cost"

# Parameter block (E/F columns)
$ws.Range("E3").Value = "Scale"
$ws.Range("F3").Value = 15
$ws.Range("E4").Value = "Base"
$ws.Range("F4").Value = 1.1
$ws.Range("E5").Value = "Exp"
$ws.Range("F5").Value = 1.4

# Level numbers + Cost formula for levels 1..20 (rows 3..22)
for ($i = 0; $i -lt 20; $i++) {
    $row = 3 + $i
    $level = $i + 1
    $ws.Range("B$row").Value = $level
    $ws.Range("C$row").Formula = '=$B' + $row + '*POWER($F$4,$F$5*$B' + $row + ')*$F$3'
}

# --- Formatting: copy styles from the "Levels and Experience" sheet -------
# (same column layout: B = level, C = value)
$wsLevels.Range("B2:C22").Copy()
$ws.Range("B2:C22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Make the new sheet the active / selected tab --------------------------
[void]$ws.Activate()
[void]$ws.Range("H7").Select()
